$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A150").Value = "2023-12-09 16:17:02"
$ws.Range("B150").Value = 0.0004

$ws.Range("A151").Value = "2023-12-09 16:17:06"
$ws.Range("B151").Value = 0.0004
